# Apply weekly update to Fruta/Hortaliza pricing rows (Ciruela - Agricola del Norte S.A. de Arica)
# The underlying data rows (2-21) are re-shuffled/updated in place: each row's
# Fecha (D), Variedad (K), Calidad (L), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P), Unidad de comercializacion (Q),
# Origen (R) and Precio $/Kg (S) are replaced with the updated weekly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 6
$ws.Cells.Item(2, 4).Value = 44243
$ws.Cells.Item(2, 11).Value = 'Black Amber'
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 300
$ws.Cells.Item(2, 14).Value = 14000
$ws.Cells.Item(2, 15).Value = 15000
$ws.Cells.Item(2, 16).Value = 14500
$ws.Cells.Item(2, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(2, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(2, 19).Value = 806

# Row 3 <- original row 10
$ws.Cells.Item(3, 4).Value = 44580
$ws.Cells.Item(3, 11).Value = 'Black Amber'
$ws.Cells.Item(3, 12).Value = 'Segunda'
$ws.Cells.Item(3, 13).Value = 270
$ws.Cells.Item(3, 14).Value = 19000
$ws.Cells.Item(3, 15).Value = 20000
$ws.Cells.Item(3, 16).Value = 19500
$ws.Cells.Item(3, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(3, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(3, 19).Value = 1083

# Row 4 <- original row 21
$ws.Cells.Item(4, 4).Value = 44278
$ws.Cells.Item(4, 11).Value = 'Angeleno'
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 300
$ws.Cells.Item(4, 14).Value = 15000
$ws.Cells.Item(4, 15).Value = 16000
$ws.Cells.Item(4, 16).Value = 15500
$ws.Cells.Item(4, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(4, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(4, 19).Value = 861

# Row 5 <- original row 3
$ws.Cells.Item(5, 4).Value = 44238
$ws.Cells.Item(5, 11).Value = 'Black Amber'
$ws.Cells.Item(5, 12).Value = 'Segunda'
$ws.Cells.Item(5, 13).Value = 300
$ws.Cells.Item(5, 14).Value = 14000
$ws.Cells.Item(5, 15).Value = 15000
$ws.Cells.Item(5, 16).Value = 14500
$ws.Cells.Item(5, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(5, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(5, 19).Value = 806

# Row 6 <- original row 4
$ws.Cells.Item(6, 4).Value = 44238
$ws.Cells.Item(6, 11).Value = 'Fortuna'
$ws.Cells.Item(6, 12).Value = 'Segunda'
$ws.Cells.Item(6, 13).Value = 300
$ws.Cells.Item(6, 14).Value = 14000
$ws.Cells.Item(6, 15).Value = 15000
$ws.Cells.Item(6, 16).Value = 14500
$ws.Cells.Item(6, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(6, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(6, 19).Value = 806

# Row 7 <- original row 11
$ws.Cells.Item(7, 4).Value = 44574
$ws.Cells.Item(7, 11).Value = 'Black Amber'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 300
$ws.Cells.Item(7, 14).Value = 18000
$ws.Cells.Item(7, 15).Value = 19000
$ws.Cells.Item(7, 16).Value = 18500
$ws.Cells.Item(7, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(7, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(7, 19).Value = 1028

# Row 8 <- original row 7
$ws.Cells.Item(8, 4).Value = 44174
$ws.Cells.Item(8, 11).Value = 'Angeleno'
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 270
$ws.Cells.Item(8, 14).Value = 20000
$ws.Cells.Item(8, 15).Value = 21000
$ws.Cells.Item(8, 16).Value = 20500
$ws.Cells.Item(8, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(8, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(8, 19).Value = 1139

# Row 9 <- original row 18
$ws.Cells.Item(9, 4).Value = 44614
$ws.Cells.Item(9, 11).Value = 'Angeleno'
$ws.Cells.Item(9, 12).Value = 'Segunda'
$ws.Cells.Item(9, 13).Value = 250
$ws.Cells.Item(9, 14).Value = 18000
$ws.Cells.Item(9, 15).Value = 19000
$ws.Cells.Item(9, 16).Value = 18500
$ws.Cells.Item(9, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(9, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(9, 19).Value = 1028

# Row 10 <- original row 13
$ws.Cells.Item(10, 4).Value = 44229
$ws.Cells.Item(10, 11).Value = 'Fortuna'
$ws.Cells.Item(10, 12).Value = 'Segunda'
$ws.Cells.Item(10, 13).Value = 300
$ws.Cells.Item(10, 14).Value = 14000
$ws.Cells.Item(10, 15).Value = 15000
$ws.Cells.Item(10, 16).Value = 14500
$ws.Cells.Item(10, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(10, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(10, 19).Value = 806

# Row 11 <- original row 2
$ws.Cells.Item(11, 4).Value = 44314
$ws.Cells.Item(11, 11).Value = 'Angeleno'
$ws.Cells.Item(11, 12).Value = 'Segunda'
$ws.Cells.Item(11, 13).Value = 250
$ws.Cells.Item(11, 14).Value = 14000
$ws.Cells.Item(11, 15).Value = 15000
$ws.Cells.Item(11, 16).Value = 14500
$ws.Cells.Item(11, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(11, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(11, 19).Value = 806

# Row 12 <- original row 9
$ws.Cells.Item(12, 4).Value = 44587
$ws.Cells.Item(12, 11).Value = 'Black Amber'
$ws.Cells.Item(12, 12).Value = 'Segunda'
$ws.Cells.Item(12, 13).Value = 300
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 16000
$ws.Cells.Item(12, 16).Value = 15500
$ws.Cells.Item(12, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(12, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(12, 19).Value = 861

# Row 13 <- original row 14
$ws.Cells.Item(13, 4).Value = 44175
$ws.Cells.Item(13, 11).Value = 'Angeleno'
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 21000
$ws.Cells.Item(13, 15).Value = 22000
$ws.Cells.Item(13, 16).Value = 21500
$ws.Cells.Item(13, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(13, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(13, 19).Value = 1194

# Row 14 <- original row 12
$ws.Cells.Item(14, 4).Value = 44169
$ws.Cells.Item(14, 11).Value = 'Angeleno'
$ws.Cells.Item(14, 12).Value = 'Tercera'
$ws.Cells.Item(14, 13).Value = 250
$ws.Cells.Item(14, 14).Value = 24000
$ws.Cells.Item(14, 15).Value = 25000
$ws.Cells.Item(14, 16).Value = 24500
$ws.Cells.Item(14, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(14, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(14, 19).Value = 1361

# Row 15 <- original row 19
$ws.Cells.Item(15, 4).Value = 44245
$ws.Cells.Item(15, 11).Value = 'Black Amber'
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 250
$ws.Cells.Item(15, 14).Value = 14000
$ws.Cells.Item(15, 15).Value = 15000
$ws.Cells.Item(15, 16).Value = 14500
$ws.Cells.Item(15, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(15, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(15, 19).Value = 806

# Row 16 <- original row 15
$ws.Cells.Item(16, 4).Value = 44239
$ws.Cells.Item(16, 11).Value = 'Fortuna'
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 13).Value = 300
$ws.Cells.Item(16, 14).Value = 15000
$ws.Cells.Item(16, 15).Value = 16000
$ws.Cells.Item(16, 16).Value = 15500
$ws.Cells.Item(16, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(16, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(16, 19).Value = 861

# Row 17 <- original row 5
$ws.Cells.Item(17, 4).Value = 44217
$ws.Cells.Item(17, 11).Value = 'Black Amber'
$ws.Cells.Item(17, 12).Value = 'Segunda'
$ws.Cells.Item(17, 13).Value = 300
$ws.Cells.Item(17, 14).Value = 16000
$ws.Cells.Item(17, 15).Value = 17000
$ws.Cells.Item(17, 16).Value = 16500
$ws.Cells.Item(17, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(17, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(17, 19).Value = 917

# Row 18 <- original row 20
$ws.Cells.Item(18, 4).Value = 44628
$ws.Cells.Item(18, 11).Value = 'Black Amber'
$ws.Cells.Item(18, 12).Value = 'Segunda'
$ws.Cells.Item(18, 13).Value = 270
$ws.Cells.Item(18, 14).Value = 15000
$ws.Cells.Item(18, 15).Value = 16000
$ws.Cells.Item(18, 16).Value = 15500
$ws.Cells.Item(18, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(18, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(18, 19).Value = 861

# Row 19 <- original row 8
$ws.Cells.Item(19, 4).Value = 44285
$ws.Cells.Item(19, 11).Value = 'Angeleno'
$ws.Cells.Item(19, 12).Value = 'Primera'
$ws.Cells.Item(19, 13).Value = 300
$ws.Cells.Item(19, 14).Value = 14000
$ws.Cells.Item(19, 15).Value = 15000
$ws.Cells.Item(19, 16).Value = 14500
$ws.Cells.Item(19, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(19, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(19, 19).Value = 806

# Row 20 <- original row 17
$ws.Cells.Item(20, 4).Value = 44596
$ws.Cells.Item(20, 11).Value = 'Black Amber'
$ws.Cells.Item(20, 12).Value = 'Segunda'
$ws.Cells.Item(20, 13).Value = 250
$ws.Cells.Item(20, 14).Value = 15000
$ws.Cells.Item(20, 15).Value = 16000
$ws.Cells.Item(20, 16).Value = 15500
$ws.Cells.Item(20, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(20, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(20, 19).Value = 861

# Row 21 <- original row 16
$ws.Cells.Item(21, 4).Value = 44650
$ws.Cells.Item(21, 11).Value = 'Angeleno'
$ws.Cells.Item(21, 12).Value = 'Segunda'
$ws.Cells.Item(21, 13).Value = 300
$ws.Cells.Item(21, 14).Value = 17000
$ws.Cells.Item(21, 15).Value = 18000
$ws.Cells.Item(21, 16).Value = 17500
$ws.Cells.Item(21, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(21, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(21, 19).Value = 972
